# Commit message: "filling marker info into status 8"
# The change fills in the "marker_1" (column J) value "NAT" for rows 2-7,
# which previously had no marker value in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the "marker_1" column (J) with "NAT" for rows 2 through 7 that were
# previously left blank.
$ws.Range("J2:J7").Value = "NAT"

# Also move the selection to G12 as recorded in the saved view state.
$ws.Range("G12").Select()
